$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TAB CONT 01")

# ---------------------------------------------------------------------------
# Header block edits
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "PRUEBA"
$ws.Range("D4").Value = 20
$ws.Range("D5").Value = 16
$ws.Range("F5").Value = 0
$ws.Range("D6").Value = 4
$ws.Range("F6").Value = 6
$ws.Range("B9").Value = "EQUIPO"

# ---------------------------------------------------------------------------
# Row 10 edits
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "ESTADO "
$ws.Range("H10").Value = "BI - T&S"
$ws.Range("L10").Value = "FX-PCX3721-0"
$ws.Range("N10").Value = 16
$ws.Range("R10").Value = 16

# ---------------------------------------------------------------------------
# Row 11 edits: text labels change, F11 (SD qty) moves to D11 (ED qty)
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "ESTADO "
$ws.Range("F11").Copy($ws.Range("D11"))
$ws.Range("D11").Value = 1
$ws.Range("F11").Clear()
$ws.Range("H11").Value = "BI - T&S"

# ---------------------------------------------------------------------------
# Row 12 edits
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = "ESTADO "
$ws.Range("N12").Value = 17
$ws.Range("R12").Value = 26

# ---------------------------------------------------------------------------
# Row 13 edits
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = "ESTADO "
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 16
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = 16

# ---------------------------------------------------------------------------
# Row 14: was the "VENT 02" sub-header; becomes a normal ESTADO data row
# while keeping the "PUNTOS SOBRANTES" summary cells (L14:R14).
# Borrow formatting for the newly populated A/B/D/H cells from row 13 (an
# already-correct ESTADO row) since row 14 didn't have those columns before.
# ---------------------------------------------------------------------------
$ws.Range("A13").Copy($ws.Range("A14"))
$ws.Range("B13").Copy($ws.Range("B14"))
$ws.Range("D13").Copy($ws.Range("D14"))
$ws.Range("H13").Copy($ws.Range("H14"))

$ws.Range("A14").Value = 5
$ws.Range("B14").Value = "ESTADO "
$ws.Range("D14").Value = 1
$ws.Range("H14").Value = "BI - T&S"
$ws.Range("M14").Value = 3
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 2
$ws.Range("Q14").Value = 4
$ws.Range("R14").Value = 10

# ---------------------------------------------------------------------------
# Row 15 edits
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = 6
$ws.Range("B15").Value = "ESTADO "
$ws.Range("H15").Value = "BI - T&S"

# ---------------------------------------------------------------------------
# Row 16 edits: text labels change, F16 (SD qty) moves to D16 (ED qty)
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "ESTADO "
$ws.Range("F16").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 1
$ws.Range("F16").Clear()
$ws.Range("H16").Value = "BI - T&S"

# ---------------------------------------------------------------------------
# Row 17 edits
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = 8
$ws.Range("B17").Value = "ESTADO "

# ---------------------------------------------------------------------------
# Row 18 edits
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = "ESTADO "

# ---------------------------------------------------------------------------
# New rows 19-25: same shape/formatting as row 15 (A/B/D/H only)
# ---------------------------------------------------------------------------
for ($i = 19; $i -le 25; $i++) {
    $ws.Range("A15").Copy($ws.Range("A$i"))
    $ws.Range("B15").Copy($ws.Range("B$i"))
    $ws.Range("D15").Copy($ws.Range("D$i"))
    $ws.Range("H15").Copy($ws.Range("H$i"))

    $ws.Range("A$i").Value = $i - 9
    $ws.Range("B$i").Value = "ESTADO "
    $ws.Range("D$i").Value = 1
    $ws.Range("H$i").Value = "BI - T&S"
}

# ---------------------------------------------------------------------------
# Remove the "Supervisor" sheet entirely
# ---------------------------------------------------------------------------
$wsSupervisor = $wb.Worksheets.Item("Supervisor")
$wsSupervisor.Delete()

Write-Host "edit complete"
